# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 816
$wsExhibit.Range("F9").Value = 390
$wsExhibit.Range("F15").Value = 13003
$wsExhibit.Range("F18").Value = 5345

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 816
$wsAll.Range("F10").Value = 390
$wsAll.Range("F17").Value = 13003
$wsAll.Range("F21").Value = 5345
